$d = $word.ActiveDocument

# Locate the paragraph containing the known anchor text, then remove the
# three empty paragraphs that immediately follow it (spacing discrepancy
# fix) while leaving the paragraph after them (with the horizontal rule
# drawing) untouched.

$anchorIndex = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Authored technical manuals and enhanced field support tools for midline hematology platforms\.") {
        $anchorIndex = $i
        break
    }
    $i = $i + 1
}

if ($anchorIndex -eq 0) {
    throw "Anchor paragraph not found"
}

# Build a Range spanning the three empty paragraphs right after the anchor
# and delete it in one shot so paragraph marks are removed along with the
# (empty) paragraph contents.
$startPara = $d.Paragraphs.Item($anchorIndex + 1)
$endPara = $d.Paragraphs.Item($anchorIndex + 3)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
